# Update the cryptos list with latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row number -> @{ D = "<price>"; E = "<volume>" } for cells that changed.
# B/C are only changed for row 51 (coin swapped out for a new one).

function Set-Cell($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# Row 2 - Bitcoin
Set-Cell 2 4 "68.056.01"
Set-Cell 2 5 "  -3.18%  "

# Row 3 - Ethereum
Set-Cell 3 4 "3.818.62"
Set-Cell 3 5 "  +1.29%  "

# Row 4 - TetherUSD
Set-Cell 4 4 "0.999"
Set-Cell 4 5 "  -0.02%  "

# Row 5 - BNB
Set-Cell 5 4 "597.64"

# Row 6 - Solana
Set-Cell 6 4 "175.04"
Set-Cell 6 5 "  -3.53%  "

# Row 7 - LidoStakedEther
Set-Cell 7 4 "3.817.43"
Set-Cell 7 5 "  +1.30%  "

# Row 8 - USDC
Set-Cell 8 5 "  +0.01%  "

# Row 9 - XRP
Set-Cell 9 5 "  -1.42%  "

# Row 10 - Dogecoin
Set-Cell 10 4 "0.161"
Set-Cell 10 5 "  -3.97%  "

# Row 11 - Toncoin
Set-Cell 11 5 "  -4.00%  "

# Row 12 - Cardano
Set-Cell 12 5 "  -4.11%  "

# Row 13 - Avalanche
Set-Cell 13 4 "38.12"
Set-Cell 13 5 "  -5.19%  "

# Row 14 - ShibaInu
Set-Cell 14 4 "0.0000246"
Set-Cell 14 5 "  -4.70%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-Cell 15 4 "4.451.20"
Set-Cell 15 5 "  +1.22%  "

# Row 16 - WrappedEther
Set-Cell 16 4 "3.810.47"
Set-Cell 16 5 "  +1.07%  "

# Row 17 - WrappedBTC
Set-Cell 17 4 "68.110.51"
Set-Cell 17 5 "  -3.11%  "

# Row 18 - TRON
Set-Cell 18 5 "  -4.59%  "

# Row 19 - Polkadot
Set-Cell 19 4 "7.16"
Set-Cell 19 5 "  -5.82%  "

# Row 20 - Chainlink
Set-Cell 20 4 "16.30"
Set-Cell 20 5 "  -2.46%  "

# Row 21 - BitcoinCash
Set-Cell 21 4 "490.45"
Set-Cell 21 5 "  -3.37%  "

# Row 22 - Uniswap
Set-Cell 22 5 "  -0.22%  "

# Row 23 - Polygon
Set-Cell 23 5 "  +0.67%  "

# Row 24 - Litecoin
Set-Cell 24 4 "84.89"
Set-Cell 24 5 "  -2.44%  "

# Row 25 - Fetch.AI
Set-Cell 25 5 "  -9.32%  "

# Row 26 - PEPE
Set-Cell 26 5 "  +2.34%  "

# Row 27 - InternetComputer(DFINITY)
Set-Cell 27 4 "12.33"
Set-Cell 27 5 "  -5.78%  "

# Row 28 - RenderToken
Set-Cell 28 4 "10.26"
Set-Cell 28 5 "  -9.98%  "

# Row 29 - Dai
Set-Cell 29 5 "  +0.12%  "

# Row 30 - PancakeSwap
Set-Cell 30 5 "  -0.66%  "

# Row 31 - ImmutableX
Set-Cell 31 4 "2.44"
Set-Cell 31 5 "  -2.61%  "

# Row 32 - EthereumClassic
Set-Cell 32 4 "32.79"
Set-Cell 32 5 "  +5.98%  "

# Row 33 - NEARProtocol
Set-Cell 33 5 "  -3.15%  "

# Row 34 - Hedera
Set-Cell 34 5 "  -4.71%  "

# Row 35 - FirstDigitalUSD
Set-Cell 35 4 "0.998"
Set-Cell 35 5 "  -0.08%  "

# Row 36 - Mantle
Set-Cell 36 5 "  -4.53%  "

# Row 37 - Kaspa
Set-Cell 37 5 "  -2.06%  "

# Row 38 - Filecoin
Set-Cell 38 5 "  -6.77%  "

# Row 39 - TheGraph
Set-Cell 39 5 "  -7.22%  "

# Row 40 - Bittensor
Set-Cell 40 4 "450.67"
Set-Cell 40 5 "  +2.29%  "

# Row 41 - OKB
Set-Cell 41 4 "48.97"
Set-Cell 41 5 "  -2.06%  "

# Row 42 - Stacks
Set-Cell 42 5 "  -4.17%  "

# Row 43 - dogwifhat
Set-Cell 43 5 "  -6.66%  "

# Row 44 - Cosmos
Set-Cell 44 4 "8.30"
Set-Cell 44 5 "  -4.44%  "

# Row 45 - Arweave
Set-Cell 45 4 "41.55"
Set-Cell 45 5 "  -9.64%  "

# Row 46 - Maker
Set-Cell 46 4 "2.834.08"
Set-Cell 46 5 "  -5.27%  "

# Row 47 - USDe (unchanged)

# Row 48 - Monero
Set-Cell 48 4 "139.12"
Set-Cell 48 5 "  +0.01%  "

# Row 49 - VeChain
Set-Cell 49 4 "0.0353"
Set-Cell 49 5 "  -3.86%  "

# Row 50 - InjectiveProtocol
Set-Cell 50 5 "  -5.47%  "

# Row 51 - coin replaced: ThetaToken -> EnergySwap
Set-Cell 51 2 "EnergySwap"
Set-Cell 51 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-Cell 51 4 "23.27"
Set-Cell 51 5 "  +2.86%  "
